$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns B/C (coin name / link) and D/E (price / volume) must remain plain
# text, matching the original inlineStr cell type. Temporarily force the affected
# ranges to Text format so Excel does not auto-convert numeric-looking strings (like
# "1.018" or "0.00001101") into real numbers, then restore the default "Normal" style
# afterwards so no stray formatting is left behind on the cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value2 = "29.222.13"
$ws.Range("E2").Value2 = "  -3.88%  "
$ws.Range("D3").Value2 = "1.964.67"
$ws.Range("E3").Value2 = "  -6.75%  "
$ws.Range("D4").Value2 = "1.018"
$ws.Range("E4").Value2 = "  +1.73%  "
$ws.Range("D5").Value2 = "328.14"
$ws.Range("E5").Value2 = "  -4.60%  "
$ws.Range("D6").Value2 = "1.016"
$ws.Range("E6").Value2 = "  +1.57%  "
$ws.Range("D7").Value2 = "0.4976"
$ws.Range("E7").Value2 = "  -6.39%  "
$ws.Range("D8").Value2 = "0.4202"
$ws.Range("E8").Value2 = "  -5.48%  "
$ws.Range("D9").Value2 = "53.56"
$ws.Range("E9").Value2 = "  -2.41%  "
$ws.Range("D10").Value2 = "0.08899"
$ws.Range("E10").Value2 = "  -5.57%  "
$ws.Range("E11").Value2 = "  -6.47%  "
$ws.Range("E12").Value2 = "  -7.58%  "
$ws.Range("D13").Value2 = "1.968.97"
$ws.Range("E13").Value2 = "  -5.38%  "
$ws.Range("D14").Value2 = "7.872"
$ws.Range("E14").Value2 = "  -8.40%  "
$ws.Range("D15").Value2 = "6.408"
$ws.Range("E15").Value2 = "  -7.63%  "
$ws.Range("D16").Value2 = "1.018"
$ws.Range("E16").Value2 = "  +1.65%  "
$ws.Range("D17").Value2 = "0.00001101"
$ws.Range("E17").Value2 = "  -5.04%  "
$ws.Range("D18").Value2 = "91.57"
$ws.Range("E18").Value2 = "  -10.16%  "
$ws.Range("E19").Value2 = "  +0.17%  "
$ws.Range("D20").Value2 = "19.25"
$ws.Range("D22").Value2 = "5.919"
$ws.Range("E22").Value2 = "  -6.73%  "
$ws.Range("D23").Value2 = "29.257.90"
$ws.Range("E23").Value2 = "  -3.81%  "
$ws.Range("D24").Value2 = "11.88"
$ws.Range("E24").Value2 = "  -5.46%  "
$ws.Range("D25").Value2 = "2.307"
$ws.Range("E25").Value2 = "  -0.24%  "
$ws.Range("D26").Value2 = "20.62"
$ws.Range("E26").Value2 = "  -5.99%  "
$ws.Range("D27").Value2 = "155.44"
$ws.Range("E27").Value2 = "  -4.51%  "
$ws.Range("E28").Value2 = "  -8.97%  "
$ws.Range("D29").Value2 = "2.287"
$ws.Range("E29").Value2 = "  -9.61%  "
$ws.Range("D30").Value2 = "126.64"
$ws.Range("E30").Value2 = "  -5.45%  "
$ws.Range("D31").Value2 = "1.049"
$ws.Range("E31").Value2 = "  -8.81%  "
$ws.Range("D32").Value2 = "0.09852"
$ws.Range("E32").Value2 = "  -6.72%  "
$ws.Range("D33").Value2 = "1.505"
$ws.Range("E33").Value2 = "  -10.33%  "
$ws.Range("D34").Value2 = "5.784"
$ws.Range("E34").Value2 = "  -7.71%  "
$ws.Range("E35").Value2 = "  -3.00%  "
$ws.Range("D36").Value2 = "0.02427"
$ws.Range("E36").Value2 = "  -8.48%  "
$ws.Range("D37").Value2 = "9.149"
$ws.Range("E37").Value2 = "  -10.51%  "
$ws.Range("B38").Value2 = "Hedera"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value2 = "0.06320"
$ws.Range("E38").Value2 = "  -7.21%  "
$ws.Range("B39").Value2 = "TrustWalletToken"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value2 = "1.286"
$ws.Range("E39").Value2 = "  -4.57%  "
$ws.Range("D40").Value2 = "0.6465"
$ws.Range("E40").Value2 = "  -8.20%  "
$ws.Range("D41").Value2 = "11.54"
$ws.Range("E41").Value2 = "  -9.35%  "
$ws.Range("D42").Value2 = "0.2014"
$ws.Range("E42").Value2 = "  -9.60%  "
$ws.Range("D43").Value2 = "1.016"
$ws.Range("E43").Value2 = "  +1.57%  "
$ws.Range("D44").Value2 = "0.6242"
$ws.Range("E44").Value2 = "  -9.36%  "
$ws.Range("D45").Value2 = "13.42"
$ws.Range("E45").Value2 = "  -7.54%  "
$ws.Range("D46").Value2 = "2.181"
$ws.Range("E46").Value2 = "  -7.18%  "
$ws.Range("D47").Value2 = "1.293"
$ws.Range("E47").Value2 = "  -6.45%  "
$ws.Range("D48").Value2 = "3.481"
$ws.Range("E48").Value2 = "  -4.43%  "
$ws.Range("E49").Value2 = "  -1.77%  "
$ws.Range("D50").Value2 = "0.06867"
$ws.Range("E50").Value2 = "  -5.44%  "
$ws.Range("D51").Value2 = "1.114"
$ws.Range("E51").Value2 = "  -8.96%  "

$dataRange.Style = "Normal"
